$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear all existing cell contents (keep formatting) so the shared-string table
# is rebuilt from scratch in the exact write order below.
$ws.Cells.ClearContents()

# Header row
$ws.Range("A1").Value = "nazev"
$ws.Range("B1").Value = "zkratka"
$ws.Range("C1").Value = "seminariciUcitIdno"

# Column A (nazev) written top-to-bottom first
$ws.Range("A2").Value = "Podnikové informační systémy"
$ws.Range("A3").Value = "Podnikové informační systémy"
$ws.Range("A4").Value = "Fyzikální praktikum C"
$ws.Range("A5").Value = "Diplomový seminář"
$ws.Range("A6").Value = "Regional geography of the Czech Republic"
$ws.Range("A7").Value = "Reg. geography of Northwestern Bohemia"
$ws.Range("A8").Value = "Počítačové modelování I"
$ws.Range("A9").Value = "Počítačové modelování I"
$ws.Range("A10").Value = "Programování A"
$ws.Range("A11").Value = "Podnikové informační systémy"
$ws.Range("A12").Value = "Podnikové informační systémy"
$ws.Range("A13").Value = "Identif. a hodn. ekosystémových služeb"
$ws.Range("A14").Value = "Podnikové informační systémy"
$ws.Range("A15").Value = "Podnikové informační systémy"
$ws.Range("A16").Value = "Podnikové informační systémy"
$ws.Range("A17").Value = "Podnikové informační systémy"
$ws.Range("A18").Value = "Sociální sítě"
$ws.Range("A19").Value = "Matematika I"
$ws.Range("A20").Value = "Identif. a hodn. ekosystémových služeb"
$ws.Range("A21").Value = "Základy autonomní robotiky"
$ws.Range("A22").Value = "Praktické aplikace hardwaru"
$ws.Range("A23").Value = "Scientific inquiry and reasoning"
$ws.Range("A24").Value = "Podnikové informační systémy"
$ws.Range("A25").Value = "Podnikové informační systémy"
$ws.Range("A26").Value = "Letní geografická škola"
$ws.Range("A27").Value = "Reflektivní seminář pedagogické praxe"
$ws.Range("A28").Value = "Reflektivní seminář pedagogické praxe"
$ws.Range("A29").Value = "Reflektivní seminář pedagogické praxe"
$ws.Range("A30").Value = "Softwarové inženýrství"
$ws.Range("A31").Value = "Softwarové inženýrství"
$ws.Range("A32").Value = "Introduction to MATLAB"
$ws.Range("A33").Value = "Odborná prezentace"
$ws.Range("A34").Value = "Odborná prezentace"

# Column B (zkratka) written top-to-bottom next
$ws.Range("B2").Value = "EIS"
$ws.Range("B3").Value = "KEIS"
$ws.Range("B4").Value = "K521"
$ws.Range("B5").Value = "K505"
$ws.Range("B6").Value = "E101"
$ws.Range("B7").Value = "E100"
$ws.Range("B8").Value = "K107"
$ws.Range("B9").Value = "P107"
$ws.Range("B10").Value = "K103"
$ws.Range("B11").Value = "EIS"
$ws.Range("B12").Value = "KEIS"
$ws.Range("B13").Value = "0153"
$ws.Range("B14").Value = "EIS"
$ws.Range("B15").Value = "KEIS"
$ws.Range("B16").Value = "EIS"
$ws.Range("B17").Value = "KEIS"
$ws.Range("B18").Value = "SON"
$ws.Range("B19").Value = "K106"
$ws.Range("B20").Value = "0153"
$ws.Range("B21").Value = "0182"
$ws.Range("B22").Value = "AHW"
$ws.Range("B23").Value = "PD101"
$ws.Range("B24").Value = "EIS"
$ws.Range("B25").Value = "KEIS"
$ws.Range("B26").Value = "0158"
$ws.Range("B27").Value = "KRSPP"
$ws.Range("B28").Value = "KSPP"
$ws.Range("B29").Value = "RSPP"
$ws.Range("B30").Value = "SWI"
$ws.Range("B31").Value = "KSWI"
$ws.Range("B32").Value = "ITM"
$ws.Range("B33").Value = "KOPRE"
$ws.Range("B34").Value = "OPRE"

# Column C (seminariciUcitIdno) numeric values, unchanged from original
$ws.Range("C2").Value = 14
$ws.Range("C3").Value = 14
$ws.Range("C4").Value = 302
$ws.Range("C5").Value = 306
$ws.Range("C6").Value = 313
$ws.Range("C7").Value = 313
$ws.Range("C8").Value = 612
$ws.Range("C9").Value = 612
$ws.Range("C10").Value = 612
$ws.Range("C11").Value = 1609
$ws.Range("C12").Value = 1609
$ws.Range("C13").Value = 2527
$ws.Range("C14").Value = 3457
$ws.Range("C15").Value = 3457
$ws.Range("C16").Value = 3606
$ws.Range("C17").Value = 3606
$ws.Range("C18").Value = 4190
$ws.Range("C19").Value = 4221
$ws.Range("C20").Value = 4625
$ws.Range("C21").Value = 4746
$ws.Range("C22").Value = 4746
$ws.Range("C23").Value = 4747
$ws.Range("C24").Value = 4991
$ws.Range("C25").Value = 4991
$ws.Range("C26").Value = 5886
$ws.Range("C27").Value = 8021
$ws.Range("C28").Value = 8021
$ws.Range("C29").Value = 8021
$ws.Range("C30").Value = 8093
$ws.Range("C31").Value = 8093
$ws.Range("C32").Value = 8514
$ws.Range("C33").Value = 8514
$ws.Range("C34").Value = 8514
